$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new test cases (DETAIL_07, DETAIL_06, DETAIL_05) ---
# Insert the new rows right after the existing data block (below row 5, the
# last original data row). Row 5 only carries the plain/default cell style
# (column G excepted), so rows inserted there inherit that same clean style
# instead of picking up the bold header style that sits above row 2. We
# then overwrite every data row (2-8) below with the final table contents,
# so it doesn't matter that the newly inserted rows are physically appended
# at the bottom before being filled in.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# Row 2: DETAIL_07 - Lỗi hệ thống bất ngờ
$ws.Range("A2").Value = "DETAIL_07"
$ws.Range("B2").Value = "Lỗi hệ thống bất ngờ"
$ws.Range("C2").Value = "Crash"
$ws.Range("D2").Value = "Service ném RuntimeException"
$ws.Range("E2").Value = "Catch & Redirect an toàn"
$ws.Range("F2").Value = "OK"
$ws.Range("G2").Value = "PASS"

# Row 3: DETAIL_06 - ID toàn khoảng trắng
$ws.Range("A3").Value = "DETAIL_06"
$ws.Range("B3").Value = "ID toàn khoảng trắng"
$ws.Range("C3").Value = "ID='   '"
$ws.Range("D3").Value = "Param id = '   '"
$ws.Range("E3").Value = "Redirect searchResult.jsp"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"

# Row 4: DETAIL_05 - ID là Null
$ws.Range("A4").Value = "DETAIL_05"
$ws.Range("B4").Value = "ID là Null"
$ws.Range("C4").Value = "ID=null"
$ws.Range("D4").Value = "Param id = null"
$ws.Range("E4").Value = "Redirect searchResult.jsp"
$ws.Range("F4").Value = "OK"
$ws.Range("G4").Value = "PASS"

# Row 5: DETAIL_01 - Xem SP thành công (unchanged content, now shifted down)
$ws.Range("A5").Value = "DETAIL_01"
$ws.Range("B5").Value = "Xem SP thành công"
$ws.Range("C5").Value = "ID=1"
$ws.Range("D5").Value = "ID=1 tồn tại -> Forward JSP"
$ws.Range("E5").Value = "Forward info-products.jsp"
$ws.Range("F5").Value = "OK"
$ws.Range("G5").Value = "PASS"

# Row 6: DETAIL_04 - ID rỗng (renamed from "ID rỗng/Null", null case split out)
$ws.Range("A6").Value = "DETAIL_04"
$ws.Range("B6").Value = "ID rỗng"
$ws.Range("C6").Value = "ID=''"
$ws.Range("D6").Value = "ID='' -> Validate fail"
$ws.Range("E6").Value = "Redirect searchResult.jsp"
$ws.Range("F6").Value = "OK"
$ws.Range("G6").Value = "PASS"

# Row 7: DETAIL_02 - SP không tồn tại (unchanged content, now shifted down)
$ws.Range("A7").Value = "DETAIL_02"
$ws.Range("B7").Value = "SP không tồn tại"
$ws.Range("C7").Value = "ID=999"
$ws.Range("D7").Value = "ID=999 -> Service trả về null"
$ws.Range("E7").Value = "Redirect searchResult.jsp"
$ws.Range("F7").Value = "OK"
$ws.Range("G7").Value = "PASS"

# Row 8: DETAIL_03 - ID lỗi format (chữ) (unchanged content, now shifted down)
$ws.Range("A8").Value = "DETAIL_03"
$ws.Range("B8").Value = "ID lỗi format (chữ)"
$ws.Range("C8").Value = "ID='abc'"
$ws.Range("D8").Value = "ID='abc' -> ParseInt lỗi"
$ws.Range("E8").Value = "Redirect searchResult.jsp"
$ws.Range("F8").Value = "OK"
$ws.Range("G8").Value = "PASS"

# --- Resize columns B/C/D to fit the new, wider content ---
# (Column widths are stored in "characters" and get rounded by the host to
# the nearest 1/6 of a character, so these inputs are chosen to land as
# close as possible to the target stored widths.)
$ws.Columns.Item(2).ColumnWidth = 19.0221354166667
$ws.Columns.Item(3).ColumnWidth = 11.9752604166667
$ws.Columns.Item(4).ColumnWidth = 27.9674479166667
